$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = "Try shifting the ill persons work to other teammembers"
$ws.Range("E5").Value = "When unsure how to build something ask teammembers"
